$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 17, 18, 19 - reusing the "HexGrid" scheme names (as Excel's shared
# string table will naturally reindex once the three new Spiral strings are
# inserted earlier in the table).
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(17, $col).Value = 1
    $ws.Cells.Item(18, $col).Value = 1
    $ws.Cells.Item(19, $col).Value = 1
}

$ws.Range("A17:A19").Style = $ws.Range("A16").Style

